# Atualização de bases das ligas, do dia: 20-06-2024 às 20:11
#
# Four pairs of adjacent match rows had their data (everything except the
# leading rank/index column A) swapped between rows, e.g. the match that was
# on row 71 is now on row 72 and vice versa. Swap columns B:AD between each
# pair of rows while leaving column A (the running rank number) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($row1, $row2) {
    $rangeA = $ws.Range("B$row1" + ":AD$row1")
    $rangeB = $ws.Range("B$row2" + ":AD$row2")

    $valuesA = $rangeA.Value2()
    $valuesB = $rangeB.Value2()

    $rangeA.Value2 = $valuesB
    $rangeB.Value2 = $valuesA
}

Swap-RowData 71 72
Swap-RowData 213 214
Swap-RowData 215 216
Swap-RowData 316 317
